$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A to fit the new, longer part description
$ws.Columns.Item(1).ColumnWidth = 65.85546875

# Add the new parts-list row (row 31): cover cap / screw terminal block
$ws.Range("A31").Value = "Aihasd 30 Stück Screw Terminal Leiterplatten-Anschlussklemme Block 5.08mm 2 Pin Pitch Klemme Printklemme 300V 16A"
$ws.Range("A31").Font.Color = 0

$ws.Range("B31").Value = 1

$ws.Hyperlinks.Add($ws.Range("C31"), "https://www.amazon.de/dp/B00VGGFPZW?ref=ppx_yo2ov_dt_b_fed_asin_title")
$ws.Range("C31").Style = "Hyperlink"

# Update selection to reflect where the cursor ended up after the edit
$ws.Range("C32").Select()
